$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2808873333333333
$ws.Range("H2").Value = 0.842662
$ws.Range("I2").Value = 0.5595554696739399
$ws.Range("J2").Value = 0.5595554696739399
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05601
$ws.Range("N2").Value = 0.16803
$ws.Range("O2").Value = 0.02710547761971223
$ws.Range("P2").Value = 0.02710547761971223
$ws.Range("Q2").Value = 0.01573249954
$ws.Range("R2").Value = 0.14159249586
$ws.Range("S2").Value = 0.01516701826023454
$ws.Range("T2").Value = 0.01516701826023454
$ws.Range("G3").Value = 0.2808873333333333
$ws.Range("H3").Value = 0.842662
$ws.Range("I3").Value = 0.5595554696739399
$ws.Range("J3").Value = 0.5595554696739399
$ws.Range("N3").Value = 5.594253
$ws.Range("O3").Value = 0.902427539668559
$ws.Range("P3").Value = 0.9024275396685592
$ws.Range("Q3").Value = 0.5237849357206666
$ws.Range("R3").Value = 4.714064421486
$ws.Range("S3").Value = 0.5049582658059386
$ws.Range("T3").Value = 0.5049582658059386
$ws.Range("G4").Value = 0.2808873333333333
$ws.Range("H4").Value = 0.842662
$ws.Range("I4").Value = 0.5595554696739399
$ws.Range("J4").Value = 0.5595554696739399
$ws.Range("M4").Value = 0.145611
$ws.Range("N4").Value = 0.436833
$ws.Range("O4").Value = 0.07046698271172858
$ws.Range("P4").Value = 0.07046698271172858
$ws.Range("Q4").Value = 0.040900285494
$ws.Range("R4").Value = 0.368102569446
$ws.Range("S4").Value = 0.03943018560776668
$ws.Range("T4").Value = 0.03943018560776668
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2210956666666667
$ws.Range("H5").Value = 0.663287
$ws.Range("I5").Value = 0.4404445303260602
$ws.Range("J5").Value = 0.4404445303260602
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05601
$ws.Range("N5").Value = 0.16803
$ws.Range("O5").Value = 0.02710547761971223
$ws.Range("P5").Value = 0.02710547761971223
$ws.Range("Q5").Value = 0.01238356829
$ws.Range("R5").Value = 0.11145211461
$ws.Range("S5").Value = 0.01193845935947769
$ws.Range("T5").Value = 0.01193845935947769
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2210956666666667
$ws.Range("H6").Value = 0.663287
$ws.Range("I6").Value = 0.4404445303260602
$ws.Range("J6").Value = 0.4404445303260602
$ws.Range("N6").Value = 5.594253
$ws.Range("O6").Value = 0.902427539668559
$ws.Range("P6").Value = 0.9024275396685592
$ws.Range("Q6").Value = 0.4122883655123333
$ws.Range("R6").Value = 3.710595289611
$ws.Range("S6").Value = 0.3974692738626206
$ws.Range("T6").Value = 0.3974692738626206
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.2210956666666667
$ws.Range("H7").Value = 0.663287
$ws.Range("I7").Value = 0.4404445303260602
$ws.Range("J7").Value = 0.4404445303260602
$ws.Range("M7").Value = 0.145611
$ws.Range("N7").Value = 0.436833
$ws.Range("O7").Value = 0.07046698271172858
$ws.Range("P7").Value = 0.07046698271172858
$ws.Range("Q7").Value = 0.032193961119
$ws.Range("R7").Value = 0.289745650071
$ws.Range("S7").Value = 0.0310367971039619
$ws.Range("T7").Value = 0.0310367971039619
